# daily auto push: 2026-01-16 02:28 UTC
#
# A new observation for 2026/01/16 (金, hour 7, ranking 201) was recorded.
# It sorts chronologically right after the existing 2026/01/16 03:00 row
# (row 659) and before the 2026/12/29 block, so a new row is inserted at
# row 660, pushing the previously-existing rows 660-701 down to 661-702.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 660:701 down to 661:702, leaving row 660 blank for the new entry.
$ws.Rows.Item(660).Insert()

# Column A holds plain text dates (e.g. "2026/01/16"), not real Excel dates.
# Force the cell to Text format *before* assigning the value so the
# "YYYY/MM/DD"-shaped string isn't auto-converted into a date serial
# number, then restore the default "Normal" style so no extra formatting
# is left behind on the cell.
$ws.Range("A660").NumberFormat = "@"
$ws.Range("A660").Value = "2026/01/16"
$ws.Range("A660").Style = "Normal"

$ws.Range("B660").Value = "金"
$ws.Range("C660").Value = 7
$ws.Range("D660").Value = 201
